$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of the model labels in column A (rows 2-26)
$labels = @{
    2  = "model_38_3_0"
    3  = "model_38_3_22"
    4  = "model_38_3_21"
    5  = "model_38_3_20"
    6  = "model_38_3_19"
    7  = "model_38_3_18"
    8  = "model_38_3_17"
    9  = "model_38_3_16"
    10 = "model_38_3_15"
    11 = "model_38_3_14"
    12 = "model_38_3_13"
    13 = "model_38_3_23"
    14 = "model_38_3_12"
    15 = "model_38_3_10"
    16 = "model_38_3_9"
    17 = "model_38_3_8"
    18 = "model_38_3_7"
    19 = "model_38_3_6"
    20 = "model_38_3_5"
    21 = "model_38_3_4"
    22 = "model_38_3_3"
    23 = "model_38_3_2"
    24 = "model_38_3_1"
    25 = "model_38_3_11"
    26 = "model_38_3_24"
}

# Metric values (B..Q) that every data row (2..26) now shares after retraining
$values = @(
    0.9999106709542923,
    0.9989035467838921,
    0.9999793340589799,
    0.9999971937300747,
    0.9999925431629452,
    0.00008338474307615185,
    0.001023490948502973,
    0.000005383998928268382,
    0.000001893807959661729,
    0.000003638575190614293,
    0.0005467859721176896,
    0.009131524685185483,
    1.000064966578696,
    0.009520272595895827,
    132.7840904041427,
    202.2600124216301
)

for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
